# Update "Förändrad" date column (C) from 2023-09-13 (serial 45182) to
# 2023-09-15 (serial 45184) for all data rows (2 through 17) on the
# active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 17; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value2 = 45184
    }
}
